$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("G5").Value = 3.6
$ws.Range("H5").Value = 3.7
$ws.Range("I5").Value = 1.95
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 2.25
$ws.Range("L5").Value = 2.6
$ws.Range("M5").Value = 1.01
$ws.Range("N5").Value = 10.5
$ws.Range("O5").Value = 1.24
$ws.Range("P5").Value = 3.3
$ws.Range("Q5").Value = 1.8
$ws.Range("R5").Value = 2
$ws.Range("S5").Value = 1.34
$ws.Range("T5").Value = 3.1
$ws.Range("U5").Value = 1.65
$ws.Range("V5").Value = 2
$ws.Range("W5").Value = 12
$ws.Range("X5").Value = 19
$ws.Range("Y5").Value = 13
$ws.Range("Z5").Value = 41
$ws.Range("AA5").Value = 29
$ws.Range("AB5").Value = 34
$ws.Range("AC5").Value = 12
$ws.Range("AD5").Value = 7
$ws.Range("AE5").Value = 15
$ws.Range("AF5").Value = 41
$ws.Range("AG5").Value = 201
$ws.Range("AH5").Value = 8
$ws.Range("AI5").Value = 10
$ws.Range("AJ5").Value = 8.5
$ws.Range("AK5").Value = 17
$ws.Range("AL5").Value = 15
$ws.Range("AM5").Value = 23
$ws.Range("AN5").Value = 5.5
$ws.Range("AO5").Value = 21
$ws.Range("AP5").Value = 26
$ws.Range("AQ5").Value = 67
$ws.Range("AR5").Value = 81
$ws.Range("AS5").Value = 151
$ws.Range("AT5").Value = 3
$ws.Range("AU5").Value = 7.5
$ws.Range("AV5").Value = 51
$ws.Range("AX5").Value = 4
$ws.Range("AY5").Value = 10
$ws.Range("AZ5").Value = 19
$ws.Range("BA5").Value = 34
$ws.Range("BB5").Value = 51
$ws.Range("BC5").Value = 126

# Row 10
$ws.Range("G10").Value = 4.72
$ws.Range("H10").Value = 3.7
$ws.Range("I10").Value = 1.55
$ws.Range("J10").Value = 5.17
$ws.Range("K10").Value = 2.2
$ws.Range("L10").Value = 2.22
$ws.Range("M10").Value = 1.02
$ws.Range("N10").Value = 10
$ws.Range("O10").Value = 1.2
$ws.Range("P10").Value = 4
$ws.Range("Q10").Value = 1.75
$ws.Range("R10").Value = 1.95
$ws.Range("S10").Value = 1.33
$ws.Range("T10").Value = 3
$ws.Range("U10").Value = 1.73
$ws.Range("V10").Value = 2

# Row 14
$ws.Range("G14").Value = 22
$ws.Range("I14").Value = 1.09
$ws.Range("J14").Value = 14
$ws.Range("L14").Value = 1.36
$ws.Range("P14").Value = 6.8
$ws.Range("Q14").Value = 1.25
$ws.Range("R14").Value = 3.6
$ws.Range("T14").Value = 4.4
$ws.Range("U14").Value = 2.05
$ws.Range("V14").Value = 1.7
$ws.Range("W14").Value = 110
$ws.Range("X14").Value = 400
$ws.Range("Y14").Value = 90
$ws.Range("AA14").Value = 450
$ws.Range("AB14").Value = 200
$ws.Range("AE14").Value = 35
$ws.Range("AG14").Value = 700
$ws.Range("AH14").Value = 13
$ws.Range("AI14").Value = 7.9
$ws.Range("AK14").Value = 6.9
$ws.Range("AM14").Value = 32
$ws.Range("AN14").Value = 20
$ws.Range("AO14").Value = 150
$ws.Range("AP14").Value = 70
$ws.Range("AT14").Value = 4.4
$ws.Range("AV14").Value = 70
$ws.Range("AY14").Value = 4.3
$ws.Range("AZ14").Value = 13.5
$ws.Range("BA14").Value = 7.9
$ws.Range("BB14").Value = 25
$ws.Range("BC14").Value = 150

# Row 15
$ws.Range("N15").Value = 10

# Row 16
$ws.Range("G16").Value = 2.3
$ws.Range("I16").Value = 3.4
$ws.Range("J16").Value = 3.2
$ws.Range("O16").Value = 1.53
$ws.Range("P16").Value = 2.38
$ws.Range("Q16").Value = 2.7
$ws.Range("R16").Value = 1.44
$ws.Range("S16").Value = 1.62
$ws.Range("T16").Value = 2.2
$ws.Range("AH16").Value = 7.5
$ws.Range("AO16").Value = 15
$ws.Range("AT16").Value = 2.2
